$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 61.2
$ws.Range("I4").Value = 61.2
$ws.Range("K4").Value = 61.2
$ws.Range("M4").Value = 52.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1256.585
$ws.Range("I40").Value = 1082.0869
$ws.Range("J40").Value = 1390.3667
$ws.Range("K40").Value = 1082.0869
$ws.Range("L40").Value = 1390.3667
$ws.Range("M40").Value = -907.0869
$ws.Range("N40").Value = -1740.3667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1564.1666
$ws.Range("I80").Value = 1350.25
$ws.Range("J80").Value = 1671.125
$ws.Range("K80").Value = 4050.75
$ws.Range("L80").Value = 5013.375
$ws.Range("M80").Value = -3052.75
$ws.Range("N80").Value = -7009.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1564.1666
$ws.Range("I83").Value = 1350.25
$ws.Range("J83").Value = 1671.125
$ws.Range("K83").Value = 12152.25
$ws.Range("L83").Value = 15040.125
$ws.Range("M83").Value = -7160.25
$ws.Range("N83").Value = -25024.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7223.722
$ws.Range("I116").Value = 8511.546
$ws.Range("J116").Value = 5200
$ws.Range("K116").Value = 8511.546
$ws.Range("L116").Value = 5200
$ws.Range("M116").Value = -5069.546
$ws.Range("N116").Value = -12084

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8799283
$ws.Range("I132").Value = 9553404
$ws.Range("K132").Value = 28660212
$ws.Range("M132").Value = -28657682

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1617.55
$ws.Range("I137").Value = 1592.1578
$ws.Range("K137").Value = 4776.4734
$ws.Range("M137").Value = -2226.4734

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3963.1562
$ws.Range("J138").Value = 4789.977
$ws.Range("L138").Value = 14369.931
$ws.Range("N138").Value = -24649.931

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 163333.33
$ws.Range("J139").Value = 250000
$ws.Range("L139").Value = 250000
$ws.Range("N139").Value = -260280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9969.462
$ws.Range("I45").Value = 15645.714
$ws.Range("K45").Value = 15645.714
$ws.Range("M45").Value = -15268.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4463.4
$ws.Range("I61").Value = 3611
$ws.Range("J61").Value = 5742
$ws.Range("K61").Value = 3611
$ws.Range("L61").Value = 5742
$ws.Range("M61").Value = -3399
$ws.Range("N61").Value = -6166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14072.429
$ws.Range("I74").Value = 2556
$ws.Range("J74").Value = 18679
$ws.Range("K74").Value = 2556
$ws.Range("L74").Value = 18679
$ws.Range("M74").Value = -1682
$ws.Range("N74").Value = -20427

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 14072.429
$ws.Range("I77").Value = 2556
$ws.Range("J77").Value = 18679
$ws.Range("K77").Value = 12780
$ws.Range("L77").Value = 93395
$ws.Range("M77").Value = -8412
$ws.Range("N77").Value = -102131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 24594.143
$ws.Range("I112").Value = 22500
$ws.Range("J112").Value = 25431.8
$ws.Range("K112").Value = 22500
$ws.Range("L112").Value = 25431.8
$ws.Range("M112").Value = -21023
$ws.Range("N112").Value = -28385.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2150.9033
$ws.Range("I132").Value = 1836.6086
$ws.Range("K132").Value = 5509.825800000001
$ws.Range("M132").Value = -2979.825800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 247147.67
$ws.Range("J133").Value = 247147.67
$ws.Range("L133").Value = 247147.67
$ws.Range("N133").Value = -252207.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 76585.22
$ws.Range("J135").Value = 76585.22
$ws.Range("L135").Value = 76585.22
$ws.Range("N135").Value = -86725.22

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4463.4
$ws.Range("I136").Value = 3611
$ws.Range("J136").Value = 5742
$ws.Range("K136").Value = 10833
$ws.Range("L136").Value = 17226
$ws.Range("M136").Value = -8283
$ws.Range("N136").Value = -22326

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 7642.6665
$ws.Range("I97").Value = 7642.6665
$ws.Range("K97").Value = 7642.6665
$ws.Range("M97").Value = -6651.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 119363.336
$ws.Range("J31").Value = 19756.334
$ws.Range("L31").Value = 19756.334
$ws.Range("N31").Value = -20346.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 119363.336
$ws.Range("J34").Value = 19756.334
$ws.Range("L34").Value = 19756.334
$ws.Range("N34").Value = -20160.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1033.05
$ws.Range("I94").Value = 1181.1428
$ws.Range("J94").Value = 687.5
$ws.Range("K94").Value = 1181.1428
$ws.Range("L94").Value = 687.5
$ws.Range("M94").Value = -730.1428000000001
$ws.Range("N94").Value = -1589.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2963.0908
$ws.Range("I99").Value = 2789.4
$ws.Range("J99").Value = 4700
$ws.Range("K99").Value = 2789.4
$ws.Range("L99").Value = 4700
$ws.Range("M99").Value = -1291.4
$ws.Range("N99").Value = -7696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 966.5
$ws.Range("I122").Value = 894.25
$ws.Range("K122").Value = 2682.75
$ws.Range("M122").Value = -232.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2963.0908
$ws.Range("I126").Value = 2789.4
$ws.Range("J126").Value = 4700
$ws.Range("K126").Value = 8368.200000000001
$ws.Range("L126").Value = 14100
$ws.Range("M126").Value = -5898.200000000001
$ws.Range("N126").Value = -19040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3846.2888
$ws.Range("I132").Value = 3768.4285
$ws.Range("K132").Value = 11305.2855
$ws.Range("M132").Value = -8775.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 16966.969
$ws.Range("I134").Value = 10021.76
$ws.Range("K134").Value = 30065.28
$ws.Range("M134").Value = -27530.28

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 141960.67
$ws.Range("J135").Value = 141960.67
$ws.Range("L135").Value = 141960.67
$ws.Range("N135").Value = -152100.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 308.875
$ws.Range("J12").Value = 295.16666
$ws.Range("L12").Value = 885.4999799999999
$ws.Range("N12").Value = -1231.49998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 10950
$ws.Range("J105").Value = 10950
$ws.Range("L105").Value = 32850
$ws.Range("N105").Value = -38092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1503.375
$ws.Range("J107").Value = 1503.375
$ws.Range("L107").Value = 4510.125
$ws.Range("N107").Value = -8350.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1547.6666
$ws.Range("J132").Value = 2055
$ws.Range("L132").Value = 18495
$ws.Range("N132").Value = -23555

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 169766.58
$ws.Range("J42").Value = 157673.2
$ws.Range("L42").Value = 157673.2
$ws.Range("N42").Value = -158643.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5998.3335
$ws.Range("I102").Value = 5247.5
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 5247.5
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -3625.5
$ws.Range("N102").Value = -10744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2920.95
$ws.Range("I113").Value = 1852.8334
$ws.Range("J113").Value = 4523.125
$ws.Range("K113").Value = 1852.8334
$ws.Range("L113").Value = 4523.125
$ws.Range("M113").Value = 317.1666
$ws.Range("N113").Value = -8863.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H115").Value = 169766.58
$ws.Range("J115").Value = 157673.2
$ws.Range("L115").Value = 157673.2
$ws.Range("N115").Value = -160023.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 19126.096
$ws.Range("I126").Value = 25945.285
$ws.Range("K126").Value = 77835.855
$ws.Range("M126").Value = -75365.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 373178.66
$ws.Range("I132").Value = 402653
$ws.Range("K132").Value = 1207959
$ws.Range("M132").Value = -1205429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 21745148
$ws.Range("I16").Value = 83335480
$ws.Range("J16").Value = 7385.5884
$ws.Range("K16").Value = 83335480
$ws.Range("L16").Value = 7385.5884
$ws.Range("M16").Value = -83335310
$ws.Range("N16").Value = -7725.5884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 430.83334
$ws.Range("I55").Value = 497.2
$ws.Range("K55").Value = 497.2
$ws.Range("M55").Value = -324.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3553.5518
$ws.Range("I61").Value = 3366.7827
$ws.Range("J61").Value = 4269.5
$ws.Range("K61").Value = 3366.7827
$ws.Range("L61").Value = 4269.5
$ws.Range("M61").Value = -3164.7827
$ws.Range("N61").Value = -4673.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3553.5518
$ws.Range("I113").Value = 3366.7827
$ws.Range("J113").Value = 4269.5
$ws.Range("K113").Value = 3366.7827
$ws.Range("L113").Value = 4269.5
$ws.Range("M113").Value = -1196.7827
$ws.Range("N113").Value = -8609.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5543.8667
$ws.Range("I132").Value = 4683.778
$ws.Range("J132").Value = 6834
$ws.Range("K132").Value = 14051.334
$ws.Range("L132").Value = 20502
$ws.Range("M132").Value = -11521.334
$ws.Range("N132").Value = -25562

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 89950
$ws.Range("J141").Value = 89950
$ws.Range("L141").Value = 89950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1574.8572
$ws.Range("I100").Value = 1482.7646
$ws.Range("K100").Value = 2965.5292
$ws.Range("M100").Value = -2424.5292

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2938.7273
$ws.Range("I132").Value = 2950.7932
$ws.Range("K132").Value = 8852.3796
$ws.Range("M132").Value = -6322.3796

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 115799.75
$ws.Range("J137").Value = 115799.75
$ws.Range("L137").Value = 115799.75
